$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C4").Value = 1430
$ws.Range("F4").Value = 1950

$ws.Range("C5").Value = 1580
$ws.Range("D5").Value = 1570
$ws.Range("G5").Value = 2000

$ws.Range("C6").Value = 1400
$ws.Range("D6").Value = 1680
$ws.Range("G6").Value = 2150

$ws.Range("D8").Value = 1200
$ws.Range("G8").Value = 730

$ws.Range("C9").Value = 1430
$ws.Range("D9").Value = 1230
$ws.Range("E9").Value = 1830
$ws.Range("F9").Value = 900

$ws.Range("G11").Select()
